$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 2 (Highs-BigM (100,100) / FEASIBLE_POINT / SOS1)
$ws.Range("E2").Value = 0.00026667
$ws.Range("F2").Value = 0.020411067
$ws.Range("G2").Value = 0.0004206910923

# Row 3 (Highs-BigM (100,100) / FEASIBLE_POINT / SOS1)
$ws.Range("E3").Value = 0.002217013
$ws.Range("F3").Value = 0.018746174
$ws.Range("G3").Value = 0.00354871923541963

# Row 4 (Highs-BigM (100,100) / FEASIBLE_POINT / Product_Mode)
$ws.Range("E4").Value = 0.006076533
$ws.Range("F4").Value = 0.029011609
$ws.Range("G4").Value = 0.007394484687407407
